# Update schedule: add Lecture 9 Data Augmentation reading/lab details,
# rework WEEK 6 lecture titles, move Midterm I, fix various lab labels,
# and record the prior ("Old") lesson titles in a new column H.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New header for the "Old" reference column (H)
$ws.Range("H1").Value = "Old"

# Lecture 9 / Data Augmentation reading link
$ws.Range("E15").Value = "09-Data-augmentation"

# Lab 5 description updated
$ws.Range("C16").Value = "Lab 5:  Gibbs, DA and Adaptive Metropolis "

# WEEK 6 lecture row: title swapped to Missing Data; prior title preserved in H
$ws.Range("C17").Value = "Lecture 10: Missing Data"
$ws.Range("H17").Value = "Lecture 10: Basics of Hypothesis Testing"

# Preserve prior lecture 11 title in H
$ws.Range("H18").Value = "Lecture 11: Hypothesis Testing"

# hw-06 now lines up with WEEK 7 / Oct 10 row; preserve prior Lec 12 title
$ws.Range("G20").Value = "hw-06"
$ws.Range("H20").Value = "Lec 12: Multiple Testing and Hierachical Models"

# Preserve prior Lec 13 title
$ws.Range("H21").Value = "Lec 13: Bayesian Multiple Testing and Hierachical Models"

# Midterm I slot becomes a review lab
$ws.Range("C22").Value = "Lab: Review"

# Midterm moved down to the Oct 19 row; hw-06 removed from here
$ws.Range("C24").Value = "Midterm 1"
$ws.Range("G24").ClearContents()

# Lab: Q&A slot becomes Lab 7: Variable Selection with hw-07
$ws.Range("C25").Value = "Lab 7: Variable Selection"
$ws.Range("G25").Value = "hw-07"

# Preserve prior Lec 14 / Lec 15 titles
$ws.Range("H26").Value = "Lec 14: Bayesian Linear Regression"
$ws.Range("H27").Value = "Lec 15: Priors in Bayesian Linear Regression"

# hw-07 removed from its old WEEK 10 row; preserve prior Lec 16 / Lec 17 titles
$ws.Range("G29").ClearContents()
$ws.Range("H29").Value = "Lec 16: Bayesian Variable Selection and Model Averaging"
$ws.Range("H30").Value = "Lec 17: Bayesian Variable Selection and Model Averaging"

# Lab: Q&A with HW 7 trimmed to Lab: Q&A
$ws.Range("C31").Value = "Lab: Q&A "

# Lec 18 / Lec 19 titles reworked; prior titles preserved in H
$ws.Range("C32").Value = "Lec 18:"
$ws.Range("H32").Value = "Lec 18: Outliers"

$ws.Range("C33").Value = "Lec 19:  Outliers"
$ws.Range("H33").Value = "Lec 19: Missing Data"

# Match the author's final selection state
$ws.Range("C22").Select()
